# dados/salas_preferenciais_2023.2.xlsx
# Melhora bastante o arquivo de saida
#
# The "CIÊNCIA DA COMPUTAÇÃO" row had its preferred-rooms list written
# without a space before each comma ("305-B, 308-B, ..."), unlike the
# corrected layout that separates each room with " , ". Fix cell B4
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "305-B , 308-B , 309-B , 310-B"

# Reflect the view state that results from having just edited/selected
# this cell: the window is scrolled one column to the right (column B
# becomes the left-most visible column) and B4 is the active selection.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B4").Select()
